$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = -0.0215
    "E" = -0.0186
    "G" = 0.3686666666666667
    "H" = 0.3686666666666667
    "I" = 0.284
    "J" = 0.2681780885780886
    "K" = 4.05
    "L" = 0.27
    "M" = 3.51
    "N" = 0.09385026737967914
    "O" = 0.8666666666666667
    "P" = 3.51
    "Q" = 0.09385026737967914
    "R" = 0.8666666666666667
    "U" = 0
    "V" = 0
    "W" = 0.1416083916083916
    "X" = 0.08725423634401565
    "Y" = 0.05435415526437594
    "Z" = 0.8571428571428571
    "AA" = 0.229866933066933
    "AB" = 0.08725423634401565
    "AC" = 0.1426126967229174
    "AE" = 0
    "AF" = 0
    "AG" = 0
    "AH" = 0
    "AI" = 0
    "AJ" = 0
    "AK" = 0
    "AP" = 0
}

foreach ($row in 2,3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
